# Append a new data row (row 13) to the sheet, mirroring the existing rows'
# formatting by copying row 12 down, then filling in the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 12
$row = 13

$ws.Range("A$srcRow:N$srcRow").Copy()
$ws.Range("A${row}:N${row}").PasteSpecial()

$ws.Cells.Item($row, 1).Value = Get-Date -Year 2016 -Month 9 -Day 6 -Hour 21 -Minute 23 -Second 23
$ws.Cells.Item($row, 2).Value = -4
$ws.Cells.Item($row, 3).Value = 53
$ws.Cells.Item($row, 4).Value = 46
$ws.Cells.Item($row, 5).Value = 53
$ws.Cells.Item($row, 6).Value = 71
$ws.Cells.Item($row, 7).Value = 13496
$ws.Cells.Item($row, 8).Value = 10671
$ws.Cells.Item($row, 9).Value = 1691
$ws.Cells.Item($row, 10).Value = 185
$ws.Cells.Item($row, 11).Value = 159
$ws.Cells.Item($row, 12).Value = 4
$ws.Cells.Item($row, 13).Value = 10
$ws.Cells.Item($row, 14).Value = "Noun"
